$d = $word.ActiveDocument

# --- Change 1: add <w:noProof/> to the run containing the inline drawing ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.InlineShapes.Count -gt 0) {
        $pp.Range.Font.NoProofing = 1
    }
}

# --- Change 2: insert new "Dangling else:" paragraph + blank paragraph ---
# Locate the paragraph whose entire text is a single space " ", right before
# the trailing blank paragraphs at the end of the document (w14:paraId="52053C98").
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $ptext = $pp.Range.Text.TrimEnd([char]13)
    if ($ptext -eq " ") {
        $targetIdx = $i
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$rng = $target.Range
$rng.Collapse(1)
$rng.InsertBefore("`r`r")

# Paragraph $targetIdx now holds the first freshly-inserted blank paragraph,
# paragraph $targetIdx+1 is the second freshly-inserted blank paragraph, and
# the original " " paragraph has shifted to $targetIdx+2.

# Insert the bold lead-in "Dangling else:" into the first new paragraph.
$p1 = $d.Paragraphs.Item($targetIdx)
$r1 = $p1.Range
$r1.Collapse(1)
$r1.InsertBefore("Dangling else:")
$r1.Font.Bold = 1
$r1.Font.BoldBi = 1

# Insert the remainder of the sentence into the second new (still separate)
# paragraph so it never touches the bold run while being formatted.
$p2 = $d.Paragraphs.Item($targetIdx + 1)
$r2 = $p2.Range
$r2.Collapse(1)
$r2.InsertBefore(" when nested if statements are used and else statement is used without proper braces then the else statement is associated with the nearest if statement and doesn’t depend on the indentation.")

# Merge the two paragraphs into one by deleting the paragraph mark between
# them (this keeps each run's own distinct formatting intact).
$p1again = $d.Paragraphs.Item($targetIdx)
$markPos = $p1again.Range.End - 1
$markRange = $d.Range($markPos, $markPos + 1)
$markRange.Delete()
